$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("df_filters")

# --- Apply the shaded "data table" style (same as the rest of the G:K
# columns) to the new block of rows before writing values into it ---
$ws.Range("G104:K113").Copy()
$ws.Range("G115:K124").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- New section: {{ df2 | rowslice(0, 3) }} ---
$ws.Range("A116").Value = "{{ df2 | rowslice(0, 3) }}"

# --- New section: {{ df2 | colslice(3) | rowslice(0, 2) }} ---
$ws.Range("A122").Value = "{{ df2 | colslice(3) | rowslice(0, 2) }}"

# --- Fill in the two previously-blank template slots that now render a
# blank-but-present ("None") column as a literal space / double-space ---
$ws.Range("I20").Value = " "
$ws.Range("J20").Value = "  "
$ws.Range("H25").Value = " "

$ws.Range("H116").Value = "name"
$ws.Range("I116").Value = "b"
$ws.Range("J116").Value = "c"
$ws.Range("K116").Value = "d"

$ws.Range("G117").Value = 0
$ws.Range("H117").Value = "a"
$ws.Range("I117").Value = 4
$ws.Range("J117").Value = 1
$ws.Range("K117").Value = 1

$ws.Range("G118").Value = 1
$ws.Range("H118").Value = "b"
$ws.Range("I118").Value = 2
$ws.Range("J118").Value = 2
$ws.Range("K118").Value = 1

$ws.Range("G119").Value = 2
$ws.Range("H119").Value = "c"
$ws.Range("I119").Value = 6
$ws.Range("J119").Value = 5
$ws.Range("K119").Value = 1

# Rows 120-121 stay blank (already shaded above).

$ws.Range("H122").Value = "d"

$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 1

$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 1

# --- View state: scroll back to the top and select G16 (matches the
# author re-reviewing the top of the sheet after adding the new filters) ---
$ws.Activate()
$ws.Range("G16").Select()
